# Adding Multiple Product Test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Update the runMode for row 3 from "Y" to "N"
$ws.Range("D3").Value = "N"

# Update the test data value in C4 (existing value "f17890k" -> "f1234567890k")
$ws.Range("C4").Value = "f1234567890k"

# Update the active selection on the LoginData sheet to B9
$ws.Activate()
$ws.Range("B9").Select()
